$wb = $excel.ActiveWorkbook

# --- Estimated sheet ---
$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Mon Dec 11 18:09:41 EST 2023"
$ws.Range("B3").Value = "Mon Dec 11 18:10:27 EST 2023"
$ws.Range("B4").Value = "Mon Dec 11 18:11:09 EST 2023"
$ws.Range("B5").Value = "Mon Dec 11 18:11:55 EST 2023"
$ws.Range("A6").Value = "Fail"
$ws.Range("B6").Value = "Mon Dec 11 18:12:38 EST 2023"
$ws.Range("A7").Value = "Fail"
$ws.Range("B7").Value = "Mon Dec 11 18:13:36 EST 2023"

# --- Existing sheet ---
$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Mon Dec 11 15:32:54 EST 2023"
$ws.Range("B3").Value = "Mon Dec 11 15:33:36 EST 2023"
$ws.Range("B4").Value = "Mon Dec 11 15:34:22 EST 2023"
$ws.Range("B5").Value = "Mon Dec 11 15:35:04 EST 2023"
$ws.Range("B6").Value = "Mon Dec 11 15:35:47 EST 2023"
$ws.Range("B7").Value = "Mon Dec 11 15:36:28 EST 2023"
$ws.Range("B8").Value = "Mon Dec 11 15:37:10 EST 2023"
$ws.Range("B9").Value = "Mon Dec 11 15:37:53 EST 2023"
$ws.Range("B10").Value = "Mon Dec 11 15:38:35 EST 2023"
$ws.Range("B11").Value = "Mon Dec 11 15:39:18 EST 2023"
$ws.Range("B12").Value = "Mon Dec 11 15:40:01 EST 2023"

# --- NewTaxReturn sheet ---
$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Mon Dec 11 15:41:21 EST 2023"
$ws.Range("B3").Value = "Mon Dec 11 15:42:03 EST 2023"
$ws.Range("B4").Value = "Mon Dec 11 15:42:47 EST 2023"
$ws.Range("B5").Value = "Mon Dec 11 15:43:28 EST 2023"
$ws.Range("B6").Value = "Mon Dec 11 15:44:10 EST 2023"
$ws.Range("B7").Value = "Mon Dec 11 15:44:51 EST 2023"
$ws.Range("B8").Value = "Mon Dec 11 15:45:32 EST 2023"
$ws.Range("B9").Value = "Mon Dec 11 15:46:12 EST 2023"
$ws.Range("B10").Value = "Mon Dec 11 15:46:53 EST 2023"
$ws.Range("B11").Value = "Mon Dec 11 15:47:34 EST 2023"
$ws.Range("B12").Value = "Mon Dec 11 15:48:19 EST 2023"
$ws.Range("B13").Value = "Mon Dec 11 15:49:01 EST 2023"
$ws.Range("B14").Value = "Mon Dec 11 15:49:43 EST 2023"
$ws.Range("B15").Value = "Mon Dec 11 15:50:24 EST 2023"
$ws.Range("B16").Value = "Mon Dec 11 15:51:05 EST 2023"

# --- Personal_IND sheet ---
$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Mon Dec 11 21:13:25 EST 2023"
$ws.Range("B4").Value = "Mon Dec 11 21:14:09 EST 2023"
$ws.Range("B5").Value = "Mon Dec 11 21:14:49 EST 2023"
$ws.Range("B6").Value = "Mon Dec 11 21:15:28 EST 2023"

# --- Personal_JNT sheet ---
$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Mon Dec 11 21:23:29 EST 2023"
$ws.Range("B4").Value = "Mon Dec 11 21:24:18 EST 2023"
$ws.Range("B5").Value = "Mon Dec 11 21:25:03 EST 2023"
$ws.Range("B6").Value = "Mon Dec 11 21:25:48 EST 2023"

# --- Personal_EL sheet ---
$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Mon Dec 11 21:12:12 EST 2023"
